$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.69394
$ws.Range("B3").Value = 7.5464799999999999
$ws.Range("B4").Value = 22.759
$ws.Range("B5").Value = 9.6023000000000014
$ws.Range("B6").Value = 6.7426400000000024
$ws.Range("B7").Value = 4.3879099999999998
$ws.Range("B8").Value = 2.5257000000000001
$ws.Range("B9").Value = 1.9352
